$wb = $excel.ActiveWorkbook

# ALC row 43: Growing Is Knowing / Growth Formula Gamma
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 21300.2
$ws.Range("I43").Value = 11999.5
$ws.Range("K43").Value = 11999.5
$ws.Range("M43").Value = -11930.5

# ALC row 76: Warding Off Temptation / Enchanted Hardsilver Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4763.9
$ws.Range("I76").Value = 4540
$ws.Range("K76").Value = 4540
$ws.Range("M76").Value = -4225

# ALC row 79: The Garden of Arcane Delights (L) / Enchanted Hardsilver Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4763.9
$ws.Range("I79").Value = 4540
$ws.Range("K79").Value = 4540
$ws.Range("M79").Value = -3448

# ALC row 88: The Grave of Hemlock Groves / Growth Formula Zeta
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 2451.5557
$ws.Range("I88").Value = 1184
$ws.Range("J88").Value = 3465.6
$ws.Range("K88").Value = 1184
$ws.Range("L88").Value = 3465.6
$ws.Range("M88").Value = -778
$ws.Range("N88").Value = -4277.6

# ALC row 91: Dappling the Highlands (L) / Growth Formula Zeta
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 2451.5557
$ws.Range("I91").Value = 1184
$ws.Range("J91").Value = 3465.6
$ws.Range("K91").Value = 1184
$ws.Range("L91").Value = 3465.6
$ws.Range("M91").Value = 220
$ws.Range("N91").Value = -6273.6

# ALC row 100: Asking for a Friend / Beetle Glue
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 7724.4
$ws.Range("I100").Value = 763
$ws.Range("K100").Value = 763
$ws.Range("M100").Value = -222

# ALC row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2753000
$ws.Range("J138").Value = 5005000
$ws.Range("L138").Value = 15015000
$ws.Range("N138").Value = -15025280

# ALC row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2816.6155
$ws.Range("I141").Value = 1783.7273
$ws.Range("K141").Value = 5351.1819
$ws.Range("M141").Value = -171.1818999999996

# ARM row 45: Hollow Hallmarks / Mythril Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2200.5
$ws.Range("I45").Value = 2029.7142
$ws.Range("J45").Value = 2499.375
$ws.Range("K45").Value = 2029.7142
$ws.Range("L45").Value = 2499.375
$ws.Range("M45").Value = -1652.7142
$ws.Range("N45").Value = -3253.375

# ARM row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2644
$ws.Range("I102").Value = 1332.3462
$ws.Range("K102").Value = 1332.3462
$ws.Range("M102").Value = 289.6538

# ARM row 120: One Foot Forward / Dwarven Mythril Shoes of Maiming
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H120").Value = 90000
$ws.Range("J120").Value = 90000
$ws.Range("L120").Value = 90000
$ws.Range("N120").Value = -99676

# BSM row 86: Through Thick and Thin / Adamantite Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1926.4286
$ws.Range("I86").Value = 1736.85
$ws.Range("J86").Value = 2400.375
$ws.Range("K86").Value = 1736.85
$ws.Range("L86").Value = 2400.375
$ws.Range("M86").Value = -613.8499999999999
$ws.Range("N86").Value = -4646.375

# BSM row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1926.4286
$ws.Range("I89").Value = 1736.85
$ws.Range("J89").Value = 2400.375
$ws.Range("K89").Value = 8684.25
$ws.Range("L89").Value = 12001.875
$ws.Range("M89").Value = -3068.25
$ws.Range("N89").Value = -23233.875

# CRP row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3397.9443
$ws.Range("I122").Value = 3413.074
$ws.Range("J122").Value = 3352.5557
$ws.Range("K122").Value = 10239.222
$ws.Range("L122").Value = 10057.6671
$ws.Range("M122").Value = -7789.222
$ws.Range("N122").Value = -14957.6671

# CUL row 131: The Mountain Steeped / Tsai tou Vounou
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1587.0952
$ws.Range("J131").Value = 1676.2632
$ws.Range("L131").Value = 5028.7896
$ws.Range("N131").Value = -15108.7896

# CUL row 137: Creative Chocolate / Gateau au Chocolat
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4982.294
$ws.Range("I137").Value = 1900
$ws.Range("J137").Value = 5174.9375
$ws.Range("K137").Value = 5700
$ws.Range("L137").Value = 15524.8125
$ws.Range("M137").Value = -600
$ws.Range("N137").Value = -25724.8125

# GSM row 40: A Little Bird Told Me / Malachite Bracelet
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 30000
$ws.Range("J40").Value = 30000
$ws.Range("L40").Value = 30000
$ws.Range("N40").Value = -30302

# GSM row 97: If I'd a Koppranickel for Every Time... / Koppranickel Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1584
$ws.Range("I97").Value = 1335.4286
$ws.Range("J97").Value = 1932
$ws.Range("K97").Value = 1335.4286
$ws.Range("L97").Value = 1932
$ws.Range("M97").Value = -839.4286
$ws.Range("N97").Value = -2924

# GSM row 102: Put the Metal to the Peddle / Durium Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 22652.709
$ws.Range("I102").Value = 1481.375
$ws.Range("J102").Value = 64995.375
$ws.Range("K102").Value = 1481.375
$ws.Range("L102").Value = 64995.375
$ws.Range("M102").Value = 140.625
$ws.Range("N102").Value = -68239.375

# GSM row 122: Awarding Academic Excellence / Ametrine
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1685.5
$ws.Range("I122").Value = 1527.3334
$ws.Range("K122").Value = 4582.0002
$ws.Range("M122").Value = -2132.0002

# GSM row 126: Gold Rush Order / Phrygian Gold Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2806.1333
$ws.Range("I126").Value = 2806.1333
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8418.3999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5948.3999
$ws.Range("N126").ClearContents()

# GSM row 132: On Board for Lar / Lar Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4109
$ws.Range("I132").Value = 3914.45
$ws.Range("K132").Value = 11743.35
$ws.Range("M132").Value = -9213.349999999999

# LTW row 7: Tan Before the Ban / Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4445.5
$ws.Range("I7").Value = 2677.5833
$ws.Range("K7").Value = 2677.5833
$ws.Range("M7").Value = -2565.5833

# LTW row 22: Skin off Their Backs / Aldgoat Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2182.682
$ws.Range("I22").Value = 604.5
$ws.Range("J22").Value = 3497.8333
$ws.Range("K22").Value = 604.5
$ws.Range("L22").Value = 3497.8333
$ws.Range("M22").Value = -309.5
$ws.Range("N22").Value = -4087.8333

# LTW row 27: Fire and Hide / Aldgoat Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2182.682
$ws.Range("I27").Value = 604.5
$ws.Range("J27").Value = 3497.8333
$ws.Range("K27").Value = 604.5
$ws.Range("L27").Value = 3497.8333
$ws.Range("M27").Value = -497.5
$ws.Range("N27").Value = -3711.8333

# LTW row 46: Supply Side Logic / Boar Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7050.6523
$ws.Range("I46").Value = 4000
$ws.Range("J46").Value = 8127.353
$ws.Range("K46").Value = 4000
$ws.Range("L46").Value = 8127.353
$ws.Range("M46").Value = -3812
$ws.Range("N46").Value = -8503.352999999999

# LTW row 82: Trainin' the Neck / Dragon Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 5649.5
$ws.Range("I82").Value = 6585.222
$ws.Range("K82").Value = 6585.222
$ws.Range("M82").Value = -6224.222

# LTW row 85: Training Is Only Skintight (L) / Dragon Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 5649.5
$ws.Range("I85").Value = 6585.222
$ws.Range("K85").Value = 6585.222
$ws.Range("M85").Value = -5337.222

# LTW row 108: Girding for Glory / Smilodonskin Trousers of Maiming
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 75000
$ws.Range("J108").Value = 100000
$ws.Range("L108").Value = 100000
$ws.Range("N108").Value = -107680

# LTW row 122: Hell on Leather / Gaja Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3536.037
$ws.Range("I122").Value = 2964.6843
$ws.Range("K122").Value = 8894.052899999999
$ws.Range("M122").Value = -6444.052899999999

# LTW row 126: Battered Books / Saiga Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4445.5
$ws.Range("I126").Value = 2677.5833
$ws.Range("K126").Value = 8032.749899999999
$ws.Range("M126").Value = -5562.749899999999

# LTW row 132: Tenets of Tanning / Silver Lobo Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2568.8645
$ws.Range("I132").Value = 2521.926
$ws.Range("J132").Value = 3075.8
$ws.Range("K132").Value = 7565.778
$ws.Range("L132").Value = 9227.400000000001
$ws.Range("M132").Value = -5035.778
$ws.Range("N132").Value = -14287.4

# WVR row 81: Where the Dragonflies, the Net Catches / Crawler Silk
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 54201.094
$ws.Range("J81").Value = 7583.7144
$ws.Range("L81").Value = 15167.4288
$ws.Range("N81").Value = -17289.4288

# WVR row 84: To Kill a Dragon on Nameday (L) / Crawler Silk
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 54201.094
$ws.Range("J84").Value = 7583.7144
$ws.Range("L84").Value = 75837.144
$ws.Range("N84").Value = -86445.144

# WVR row 100: Of Great Import / Kudzu Thread
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1092.0952
$ws.Range("I100").Value = 1125.5
$ws.Range("K100").Value = 2251
$ws.Range("M100").Value = -1710

# WVR row 122: Heavy Armoire / Dark Hempen Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3763.4482
$ws.Range("J122").Value = 4470.7144
$ws.Range("L122").Value = 13412.1432
$ws.Range("N122").Value = -18312.1432

# WVR row 126: A Polished Purchase / Snow Linen
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2427.0334
$ws.Range("I126").Value = 2151.5715
$ws.Range("K126").Value = 6454.7145
$ws.Range("M126").Value = -3984.7145
